$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.475.74"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.361.52"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'573.14"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'137.29"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.357.19"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "'7.49"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "'0.388"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "3.931.45"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "'26.01"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "3.359.13"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "61.506.73"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'13.97"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "'9.34"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "'379.35"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").Value = "'0.552"
$ws.Range("E23").Value = "  -3.81%  "
$ws.Range("D24").Value = "3.499.64"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").Value = "'71.25"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +10.53%  "
$ws.Range("D29").Value = "'7.48"
$ws.Range("E29").Value = "  -4.30%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").Value = "'8.18"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'23.63"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'5.21"
$ws.Range("E36").Value = "  -6.44%  "
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").Value = "'165.23"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'0.0764"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "'0.770"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("D45").Value = "'41.46"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'4.39"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "'23.83"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").Value = "'22.87"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "2.373.64"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'0.0260"
$ws.Range("E51").Value = "  -2.71%  "
